$wb = $excel.ActiveWorkbook

# --- ALC (sheet1) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 7341.643
$ws.Range("I33").Value = 42.6
$ws.Range("K33").Value = 42.6
$ws.Range("M33").Value = 186.4
$ws.Range("H40").Value = 1678.8572
$ws.Range("I40").Value = 1600.2858
$ws.Range("J40").Value = 1836
$ws.Range("K40").Value = 1600.2858
$ws.Range("L40").Value = 1836
$ws.Range("M40").Value = -1425.2858
$ws.Range("N40").Value = -2186
$ws.Range("H64").Value = 960261.6
$ws.Range("I64").Value = 1506536.4
$ws.Range("J64").Value = 4280.75
$ws.Range("K64").Value = 1506536.4
$ws.Range("L64").Value = 4280.75
$ws.Range("M64").Value = -1506288.4
$ws.Range("N64").Value = -4776.75
$ws.Range("H67").Value = 960261.6
$ws.Range("I67").Value = 1506536.4
$ws.Range("J67").Value = 4280.75
$ws.Range("K67").Value = 1506536.4
$ws.Range("L67").Value = 4280.75
$ws.Range("M67").Value = -1505678.4
$ws.Range("N67").Value = -5996.75
$ws.Range("H100").Value = 1738.0294
$ws.Range("I100").Value = 1040.3704
$ws.Range("J100").Value = 4429
$ws.Range("K100").Value = 1040.3704
$ws.Range("L100").Value = 4429
$ws.Range("M100").Value = -499.3704
$ws.Range("N100").Value = -5511
$ws.Range("H125").Value = 1491.5834
$ws.Range("I125").Value = 1440
$ws.Range("J125").Value = 1528.4286
$ws.Range("K125").Value = 12960
$ws.Range("L125").Value = 13755.8574
$ws.Range("M125").Value = -10500
$ws.Range("N125").Value = -18675.8574
$ws.Range("H132").Value = 1554.317
$ws.Range("I132").Value = 1459.9487
$ws.Range("K132").Value = 4379.8461
$ws.Range("M132").Value = -1849.8461
$ws.Range("H140").Value = 78566
$ws.Range("J140").Value = 78566
$ws.Range("L140").Value = 78566
$ws.Range("N140").Value = -88926

# --- ARM (sheet2) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2162.5417
$ws.Range("I61").Value = 2106.1052
$ws.Range("K61").Value = 2106.1052
$ws.Range("M61").Value = -1894.1052
$ws.Range("H136").Value = 2162.5417
$ws.Range("I136").Value = 2106.1052
$ws.Range("K136").Value = 6318.3156
$ws.Range("M136").Value = -3768.3156

# --- BSM (sheet3) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1833.1666
$ws.Range("I134").Value = 1866.6666
$ws.Range("J134").Value = 1799.6666
$ws.Range("K134").Value = 5599.9998
$ws.Range("L134").Value = 5398.9998
$ws.Range("M134").Value = -3064.9998
$ws.Range("N134").Value = -10468.9998

# --- CRP (sheet4) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3538.3572
$ws.Range("I31").Value = 2996.923
$ws.Range("J31").Value = 4007.6
$ws.Range("K31").Value = 2996.923
$ws.Range("L31").Value = 4007.6
$ws.Range("M31").Value = -2701.923
$ws.Range("N31").Value = -4597.6
$ws.Range("H34").Value = 3538.3572
$ws.Range("I34").Value = 2996.923
$ws.Range("J34").Value = 4007.6
$ws.Range("K34").Value = 2996.923
$ws.Range("L34").Value = 4007.6
$ws.Range("M34").Value = -2794.923
$ws.Range("N34").Value = -4411.6
$ws.Range("H58").Value = 66667916
$ws.Range("I58").Value = 125000650
$ws.Range("J58").Value = 1941.4286
$ws.Range("K58").Value = 125000650
$ws.Range("L58").Value = 1941.4286
$ws.Range("M58").Value = -125000447
$ws.Range("N58").Value = -2347.4286
$ws.Range("H94").Value = 4360.421
$ws.Range("I94").Value = 700
$ws.Range("J94").Value = 4791.0586
$ws.Range("K94").Value = 700
$ws.Range("L94").Value = 4791.0586
$ws.Range("M94").Value = -249
$ws.Range("N94").Value = -5693.0586
$ws.Range("H132").Value = 2112.0977
$ws.Range("I132").Value = 1288.5714
$ws.Range("J132").Value = 6916
$ws.Range("K132").Value = 3865.7142
$ws.Range("L132").Value = 20748
$ws.Range("M132").Value = -1335.7142
$ws.Range("N132").Value = -25808
$ws.Range("H134").Value = 2041.7142
$ws.Range("I134").Value = 2158.4
$ws.Range("J134").Value = 1750
$ws.Range("K134").Value = 6475.200000000001
$ws.Range("L134").Value = 5250
$ws.Range("M134").Value = -3940.200000000001
$ws.Range("N134").Value = -10320
$ws.Range("H136").Value = 66667916
$ws.Range("I136").Value = 125000650
$ws.Range("J136").Value = 1941.4286
$ws.Range("K136").Value = 375001950
$ws.Range("L136").Value = 5824.2858
$ws.Range("M136").Value = -374999400
$ws.Range("N136").Value = -10924.2858

# --- CUL (sheet5) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 765.55
$ws.Range("J131").Value = 823.16095
$ws.Range("L131").Value = 2469.48285
$ws.Range("N131").Value = -12549.48285

# --- GSM (sheet6) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2471
$ws.Range("I122").Value = 2529.4
$ws.Range("K122").Value = 7588.200000000001
$ws.Range("M122").Value = -5138.200000000001
$ws.Range("H140").Value = 41333.8
$ws.Range("J140").Value = 41333.8
$ws.Range("L140").Value = 41333.8
$ws.Range("N140").Value = -51693.8

# --- LTW (sheet7) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2195.3462
$ws.Range("I68").Value = 2206
$ws.Range("J68").Value = 2180.818
$ws.Range("K68").Value = 2206
$ws.Range("L68").Value = 2180.818
$ws.Range("M68").Value = -1457
$ws.Range("N68").Value = -3678.818
$ws.Range("H71").Value = 2195.3462
$ws.Range("I71").Value = 2206
$ws.Range("J71").Value = 2180.818
$ws.Range("K71").Value = 11030
$ws.Range("L71").Value = 10904.09
$ws.Range("M71").Value = -7286
$ws.Range("N71").Value = -18392.09
$ws.Range("H136").Value = 3607
$ws.Range("I136").Value = 3946.5
$ws.Range("J136").Value = 3199.6
$ws.Range("K136").Value = 11839.5
$ws.Range("L136").Value = 9598.799999999999
$ws.Range("M136").Value = -9289.5
$ws.Range("N136").Value = -14698.8
$ws.Range("H138").Value = 53606.332
$ws.Range("J138").Value = 53606.332
$ws.Range("L138").Value = 53606.332
$ws.Range("N138").Value = -63886.332
$ws.Range("H139").Value = 49875
$ws.Range("J139").Value = 49875
$ws.Range("L139").Value = 49875
$ws.Range("N139").Value = -60155

# --- WVR (sheet8) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 795.275
$ws.Range("I136").Value = 553.41174
$ws.Range("J136").Value = 2165.8333
$ws.Range("K136").Value = 1660.23522
$ws.Range("L136").Value = 6497.499899999999
$ws.Range("M136").Value = 889.76478
$ws.Range("N136").Value = -11597.4999
$ws.Range("H138").Value = 56400
$ws.Range("J138").Value = 56400
$ws.Range("L138").Value = 56400
$ws.Range("N138").Value = -66680
$ws.Range("H139").Value = 46488.89
$ws.Range("J139").Value = 46488.89
$ws.Range("L139").Value = 46488.89
$ws.Range("N139").Value = -56768.89
$ws.Range("H141").Value = 80000
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()
